$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.727.19'
$ws.Cells.Item(2, 5).Value = '  +0.07%  '

$ws.Cells.Item(3, 4).Value = '3.849.99'
$ws.Cells.Item(3, 5).Value = '  -0.27%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).Value = '''601.91'

$ws.Cells.Item(6, 4).Value = '''170.27'
$ws.Cells.Item(6, 5).Value = '  +0.88%  '

$ws.Cells.Item(7, 4).Value = '3.850.09'
$ws.Cells.Item(7, 5).Value = '  -0.33%  '

$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 5).Value = '  -0.57%  '

$ws.Cells.Item(10, 5).Value = '  +0.69%  '

$ws.Cells.Item(11, 5).Value = '  +2.22%  '

$ws.Cells.Item(12, 5).Value = '  -0.81%  '

$ws.Cells.Item(13, 4).Value = '''0.0000280'
$ws.Cells.Item(13, 5).Value = '  +11.55%  '

$ws.Cells.Item(14, 4).Value = '''36.95'
$ws.Cells.Item(14, 5).Value = '  -1.78%  '

$ws.Cells.Item(15, 4).Value = '4.492.65'

$ws.Cells.Item(16, 4).Value = '3.845.66'
$ws.Cells.Item(16, 5).Value = '  -0.07%  '

$ws.Cells.Item(17, 4).Value = '68.647.37'

$ws.Cells.Item(18, 4).Value = '''18.32'
$ws.Cells.Item(18, 5).Value = '  +0.07%  '

$ws.Cells.Item(19, 4).Value = '''7.37'
$ws.Cells.Item(19, 5).Value = '  -2.82%  '

$ws.Cells.Item(20, 5).Value = '  -0.53%  '

$ws.Cells.Item(21, 4).Value = '''10.95'
$ws.Cells.Item(21, 5).Value = '  +0.62%  '

$ws.Cells.Item(22, 4).Value = '''474.14'
$ws.Cells.Item(22, 5).Value = '  -0.31%  '

$ws.Cells.Item(23, 4).Value = '''0.726'
$ws.Cells.Item(23, 5).Value = '  -1.89%  '

$ws.Cells.Item(24, 5).Value = '  +1.50%  '

$ws.Cells.Item(25, 4).Value = '''83.50'
$ws.Cells.Item(25, 5).Value = '  -1.59%  '

$ws.Cells.Item(26, 4).Value = '''2.24'
$ws.Cells.Item(26, 5).Value = '  -0.85%  '

$ws.Cells.Item(27, 4).Value = '''12.11'
$ws.Cells.Item(27, 5).Value = '  -2.72%  '

$ws.Cells.Item(28, 4).Value = '''10.35'
$ws.Cells.Item(28, 5).Value = '  +2.25%  '

$ws.Cells.Item(29, 5).Value = '  +0.03%  '

$ws.Cells.Item(30, 5).Value = '  -0.28%  '

$ws.Cells.Item(31, 4).Value = '3.996.97'
$ws.Cells.Item(31, 5).Value = '  -0.42%  '

$ws.Cells.Item(32, 5).Value = '  -0.83%  '

$ws.Cells.Item(33, 5).Value = '  +0.62%  '

$ws.Cells.Item(34, 4).Value = '''2.29'
$ws.Cells.Item(34, 5).Value = '  -1.46%  '

$ws.Cells.Item(35, 4).Value = '''9.30'
$ws.Cells.Item(35, 5).Value = '  -2.35%  '

$ws.Cells.Item(36, 4).Value = '3.811.67'
$ws.Cells.Item(36, 5).Value = '  -0.49%  '

$ws.Cells.Item(37, 4).Value = '''3.82'
$ws.Cells.Item(37, 5).Value = '  +14.13%  '

$ws.Cells.Item(38, 5).Value = '  -1.78%  '

$ws.Cells.Item(39, 2).Value = 'Kaspa'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(39, 4).Value = '''0.140'
$ws.Cells.Item(39, 5).Value = '  -1.21%  '

$ws.Cells.Item(40, 2).Value = 'Mantle'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(40, 4).Value = '''1.02'
$ws.Cells.Item(40, 5).Value = '  +0.45%  '

$ws.Cells.Item(41, 4).Value = '''5.92'
$ws.Cells.Item(41, 5).Value = '  -1.56%  '

$ws.Cells.Item(42, 4).Value = '''1.00'
$ws.Cells.Item(42, 5).Value = '  -0.07%  '

$ws.Cells.Item(43, 4).Value = '''0.315'
$ws.Cells.Item(43, 5).Value = '  -0.84%  '

$ws.Cells.Item(44, 2).Value = 'FLOKI'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Cells.Item(44, 4).Value = '''0.000301'
$ws.Cells.Item(44, 5).Value = '  +10.15%  '

$ws.Cells.Item(45, 2).Value = 'Stacks'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(45, 4).Value = '''2.00'
$ws.Cells.Item(45, 5).Value = '  -1.30%  '

$ws.Cells.Item(46, 2).Value = 'USDe'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(46, 4).Value = '''1.00'
$ws.Cells.Item(46, 5).Value = '  +0.00%  '

$ws.Cells.Item(47, 4).Value = '''420.93'
$ws.Cells.Item(47, 5).Value = '  -2.40%  '

$ws.Cells.Item(48, 4).Value = '''8.68'
$ws.Cells.Item(48, 5).Value = '  +0.37%  '

$ws.Cells.Item(49, 4).Value = '''46.89'
$ws.Cells.Item(49, 5).Value = '  -1.63%  '

$ws.Cells.Item(50, 2).Value = 'Monero'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(50, 4).Value = '''141.77'
$ws.Cells.Item(50, 5).Value = '  -0.33%  '

$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 4).Value = '''26.07'
$ws.Cells.Item(51, 5).Value = '  +3.83%  '
